$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying comment rows got reshuffled: row 7 <-> row 8, and row 20 <-> row 21
# (columns B:D - the comment text, author, and like-count). Swap them in place,
# using an out-of-the-way scratch range so shared-string cell types (text) are
# preserved via Copy instead of Value assignment (which would coerce numeric-
# looking strings like "1"/"2"/"3" into real numbers).

$ws.Range("B7:D7").Copy($ws.Range("B200:D200"))
$ws.Range("B8:D8").Copy($ws.Range("B7:D7"))
$ws.Range("B200:D200").Copy($ws.Range("B8:D8"))

$ws.Range("B20:D20").Copy($ws.Range("B200:D200"))
$ws.Range("B21:D21").Copy($ws.Range("B20:D20"))
$ws.Range("B200:D200").Copy($ws.Range("B21:D21"))

$ws.Range("B200:D200").Clear()
